# Covid-19 Calabria - Vibo Valentia workbook update
# Adds new daily data (rows 476-484) to the 4 sheets and extends the
# date column (A) through row 510 (new days with no data yet).
#
# Sheet layout (same on all 4 sheets):
#   A = date (serial number, dd/mm/yyyy display)
#   C = daily count
#   D = 7-day rolling average = AVERAGE(C[n-6]:C[n])

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New daily values (column C) for rows 476-484, per sheet (1-based sheet
# index matches the tab order: 1 Nuovi casi, 2 Deceduti, 3 Dimessi Guariti,
# 4 Ricoveri).
# ---------------------------------------------------------------------------
$newValues = @{
    1 = @{476=3; 477=0; 478=4; 479=1; 480=0; 481=1; 482=0; 483=0; 484=1}
    2 = @{476=0; 477=0; 478=0; 479=0; 480=0; 481=0; 482=0; 483=0; 484=0}
    3 = @{476=5; 477=4; 478=13; 479=3; 480=3; 481=5; 482=0; 483=8; 484=5}
    4 = @{476=3; 477=3; 478=2; 479=3; 480=3; 481=3; 482=4; 483=4; 484=4}
}

# Base date serial for row 476 is 27/06/2021 (44374); each following row is
# one calendar day later, all the way through row 510 (21/07/2021 = 44408).
$baseRow = 476
$baseSerial = 44374
$lastDateRow = 510

for ($sheetIdx = 1; $sheetIdx -le 4; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $values = $newValues[$sheetIdx]

    # --- Column A: make sure every day through row 510 has its date -------
    for ($r = $baseRow; $r -le $lastDateRow; $r++) {
        $serial = $baseSerial + ($r - $baseRow)
        $ws.Cells.Item($r, 1).Value = $serial
    }

    # --- Columns C (value) and D (7-day rolling average formula) ----------
    foreach ($r in 476..484) {
        $ws.Cells.Item($r, 3).Value = $values[$r]

        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "#,##0"
        $dCell.Font.Color = 0
        $startRow = $r - 6
        $dCell.Formula = "=AVERAGE(C$startRow`:C$r)"
    }
EOF_PLACEHOLDER
}
